$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = "Personnel Subarea"
$ws.Range("M2").Value = "Warsaw"
$ws.Range("M3").Value = "Lodz"
$ws.Range("M4").Value = "Poznan"
